# Oppdatert til vår 2023
# Adds the "2023 - Høst" exam row below the existing "2023 - Vår" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A15").Value = "2023 - Høst"
$ws.Range("B15").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-23-h.pdf)"
$ws.Range("C15").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-23-h-fasit.pdf)"

$ws.Range("C16").Select()
